$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 1050
$ws.Range("F8").Value = 1659
$ws.Range("F9").Value = 386
$ws.Range("F10").Value = 1763
$ws.Range("F12").Value = 1215
$ws.Range("F15").Value = 2289
$ws.Range("F16").Value = 314
$ws.Range("F18").Value = 1047
$ws.Range("F19").Value = 568
$ws.Range("F21").Value = 1416
$ws.Range("F22").Value = 1174
$ws.Range("F23").Value = 126
$ws.Range("F25").Value = 1288
$ws.Range("F26").Value = 938
$ws.Range("F28").Value = 1251
$ws.Range("F29").Value = 125
$ws.Range("F30").Value = 1206
$ws.Range("F31").Value = 390
$ws.Range("F36").Value = 414
$ws.Range("F37").Value = 23
$ws.Range("F40").Value = 2160
$ws.Range("F41").Value = 115
$ws.Range("F42").Value = 867
$ws.Range("F43").Value = 1846
$ws.Range("F45").Value = 826

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 57
$ws.Range("F6").Value = 290
$ws.Range("F8").Value = 85
$ws.Range("F9").Value = 23
$ws.Range("F15").Value = 100421
$ws.Range("F19").Value = 50
$ws.Range("F20").Value = 50
$ws.Range("F21").Value = 204
$ws.Range("F22").Value = 273
$ws.Range("F24").Value = 255
$ws.Range("F26").Value = 71
$ws.Range("F27").Value = 63
$ws.Range("F31").Value = 30
$ws.Range("F32").Value = 201
$ws.Range("F36").Value = 79

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 2981
$ws.Range("F6").Value = 4785
$ws.Range("F10").Value = 861
$ws.Range("F11").Value = 514
$ws.Range("F12").Value = 528
$ws.Range("F13").Value = 1245
$ws.Range("F14").Value = 359
$ws.Range("F15").Value = 973

$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 4785
$ws.Range("F7").Value = 861
$ws.Range("F8").Value = 514
$ws.Range("F10").Value = 528
$ws.Range("F11").Value = 1245
$ws.Range("F12").Value = 1050
$ws.Range("F14").Value = 1659
$ws.Range("F15").Value = 386
$ws.Range("F16").Value = 85
$ws.Range("F17").Value = 1763
$ws.Range("F19").Value = 1215
$ws.Range("F20").Value = 23
$ws.Range("F21").Value = 973
$ws.Range("F22").Value = 973
$ws.Range("F23").Value = 2289
$ws.Range("F25").Value = 314
$ws.Range("F27").Value = 1047
$ws.Range("F28").Value = 568
$ws.Range("F29").Value = 1416
$ws.Range("F31").Value = 1174
$ws.Range("F32").Value = 126
$ws.Range("F33").Value = 1288
$ws.Range("F34").Value = 938
$ws.Range("F35").Value = 1251
$ws.Range("F36").Value = 125
$ws.Range("F37").Value = 50
$ws.Range("F38").Value = 1206
$ws.Range("F39").Value = 390
$ws.Range("F44").Value = 23
$ws.Range("F46").Value = 2160
$ws.Range("F47").Value = 115
$ws.Range("F48").Value = 867
$ws.Range("F49").Value = 1846
$ws.Range("F50").Value = 826
